function Set-DateText {
    param($range, $text)
    $parts = $text -split '-'
    $d = [int]$parts[0]
    $m = [int]$parts[1]
    if ($d -le 12 -and $m -le 12) {
        # Ambiguous DD-MM-YYYY pattern that Excel could mis-parse as a date;
        # force literal text via a leading apostrophe (quote-prefix).
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(37, 35, '22-12-2025', 'SHAFEEK', 9995493270, '04-01-2026', 'Akash R', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', $null),
    @(38, 36, '22-12-2025', 'shiraf', 8089328458, '08-01-2026', 'MOHAMMED NABEEL N', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', 'discuss with family and visit later'),
    @(39, 37, '22-12-2025', 'SHAFEEK', 9567303535, '27-12-2025', 'Akash R', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', $null),
    @(40, 38, '22-12-2025', 'Bharath', 8590980810, '26-01-2026', 'MOHAMMED NABEEL N', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', 'confirm within days'),
    @(41, 39, '23-12-2025', 'ANIL', 8590150651, '07-01-2026', 'ARJUN G.S', 'Loss', 'ENQUIRY', 'Enquiry for Relative/Friend', '-', $null),
    @(42, 40, '23-12-2025', 'SUCKIN OUSAP', 8848865920, '29-12-2025', 'MOHAMMED NABEEL N', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', $null),
    @(43, 41, '23-12-2025', 'SHAKIR', 7736376366, '29-01-2026', 'NIHAL S', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', 'VISIT WITH FAMILY'),
    @(44, 42, '24-12-2025', 'RAEES', 8089599898, '18-01-2026', 'Deepu M', 'Loss', 'PRODUCT', 'Product Already Booked', '-', $null),
    @(45, 43, '24-12-2025', 'ROSHAN', 8606000282, '01-01-2026', 'Deepu M', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', $null),
    @(46, 44, '24-12-2025', 'NABEEL', 9895727494, '15-01-2026', 'Akash R', 'Loss', 'PRODUCT', 'REQUIRED MODEL NOT AVAILABLE', '-', $null),
    @(47, 45, '25-12-2025', 'SHAD', 7356570554, '22-03-2026', 'NIHAL S', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT TRIAL', '-', 'VISIT WITH FAM LATER'),
    @(48, 46, '25-12-2025', 'ASIF', 7012392302, '11-01-2026', 'MOHAMMED NABEEL N', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', $null),
    @(49, 47, '25-12-2025', 'MIDHIN', 7034992304, '01-02-2026', 'Deepu M', 'Loss', 'ENQUIRY', 'ENQUIRY WITHOUT BRIDE/FAMILY', '-', $null)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "0"
    Set-DateText $ws.Cells.Item($r, 2) $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 4).NumberFormat = "0"
    Set-DateText $ws.Cells.Item($r, 5) $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    if ($row[11] -ne $null) {
        $ws.Cells.Item($r, 11).Value = $row[11]
    }
}
